$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D-column values that would otherwise be parsed as numbers
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

# Apply the updated cell values from the crypto data refresh
$ws.Range('D2').Value = '30.550.49'
$ws.Range('E2').Value = '  +2.07%  '
$ws.Range('D3').Value = '1.672.18'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '219.50'
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '29.54'
$ws.Range('E8').Value = '  +3.48%  '
$ws.Range('E9').Value = '  +2.11%  '
$ws.Range('D10').Value = '0.0634'
$ws.Range('E10').Value = '  +4.29%  '
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').Value = '1.913.27'
$ws.Range('E12').Value = '  +2.46%  '
$ws.Range('D13').Value = '1.682.02'
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('D14').Value = '0.613'
$ws.Range('E14').Value = '  +8.80%  '
$ws.Range('D15').Value = '10.13'
$ws.Range('E15').Value = '  +9.05%  '
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').Value = '30.558.38'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').Value = '66.10'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('D19').Value = '243.07'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = '0.0₃0722'
$ws.Range('E20').Value = '  +2.92%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '4.24'
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('D23').Value = '9.96'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').Value = '157.80'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').Value = '15.85'
$ws.Range('E26').Value = '  +2.12%  '
$ws.Range('D28').Value = '6.66'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.0493'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('E31').Value = '  +2.57%  '
$ws.Range('E32').Value = '  +2.66%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '3.28'
$ws.Range('E33').Value = '  +3.46%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.496.49'
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('E35').Value = '  +6.93%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').Value = '83.59'
$ws.Range('E37').Value = '  +10.40%  '
$ws.Range('E38').Value = '  +4.92%  '
$ws.Range('E39').Value = '  +7.07%  '
$ws.Range('E40').Value = '  -3.76%  '
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('E42').Value = '  +1.03%  '
$ws.Range('E43').Value = '  +1.91%  '
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('E47').Value = '  +3.79%  '
$ws.Range('D48').Value = '51.24'
$ws.Range('E48').Value = '  -3.02%  '
$ws.Range('D49').Value = '1.804.74'
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('D50').Value = '94.40'
$ws.Range('E50').Value = '  +5.27%  '
$ws.Range('E51').Value = '  -0.56%  '
